$d = $word.ActiveDocument

# 1. Remove "HTML, CSS, " from the "Python, R (R Shiny), SQL / MySQL, HTML, CSS, Shell scripting" line
$d.Content.Find.Execute("HTML, CSS, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Insert "Flask, " before "AWS" in "PyTorch, TensorFlow, Keras, Hadoop, AWS, Google Analytics, Kubernetes"
$d.Content.Find.Execute("Hadoop, AWS", $false, $false, $false, $false, $false, $true, 1, $false, "Hadoop, Flask, AWS", 2)
